# Add 2022-Q4 data: a new per-quarter sheet "2022-Q4" plus a new summary
# row at the top of the "总计" (totals) sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" sheet: insert a new row 2 for "2022-Q4", shifting the
#    existing quarter rows (old rows 2..9) down to rows 3..10.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# Shift B:D down one row at a time, bottom-up so we never clobber a row
# before it has been copied.
for ($r = 9; $r -ge 2; $r--) {
    $total.Range("B" + $r + ":D" + $r).Copy($total.Range("B" + ($r + 1)))
}

# Rebuild the A column (0-based row index) for rows 2..10, reusing the
# existing index-cell style.
$total.Range("A2").Copy($total.Range("A3:A10"))
for ($r = 2; $r -le 10; $r++) {
    $total.Cells.Item($r, 1).Value = $r - 2
}

# New top row: 2022-Q4 summary.
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 5
$total.Cells.Item(2, 4).Value = 0.06

# ---------------------------------------------------------------------
# 2) New "2022-Q4" worksheet: duplicate the "2022-Q3" sheet (same
#    column layout/formatting) right after "总计", rename it, then
#    overwrite its data rows with the 2022-Q4 fund holdings.
# ---------------------------------------------------------------------
$srcQuarter = $wb.Worksheets.Item("2022-Q3")
$srcQuarter.Copy($null, $total)

$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

$q4Rows = @(
    @("161715", "招商中证大宗商品股票指数（LOF）", "1.69", "94.56", "1.11", "0.0188", 1),
    @("000646", "华润元大量化优选混合A", "0.19", "65.16", "8.49", "0.0161", 3),
    @("540004", "汇丰晋信2026周期混合", "1.08", "23.97", "1.16", "0.0125", 8),
    @("159990", "银华巨潮小盘价值ETF", "0.78", "97.02", "1.08", "0.0084", 8),
    @("007827", "华润元大量化优选混合C", "0.01", "65.16", "8.49", "0.0008", 3)
)

$r = 2
foreach ($row in $q4Rows) {
    $q4.Cells.Item($r, 2).Value = $row[0]
    $q4.Cells.Item($r, 3).Value = $row[1]
    $q4.Cells.Item($r, 4).Value = $row[2]
    $q4.Cells.Item($r, 5).Value = $row[3]
    $q4.Cells.Item($r, 6).Value = $row[4]
    $q4.Cells.Item($r, 7).Value = $row[5]
    $q4.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}
